# issue #5: add legislator_id, name, date into dataframe
# The "股票" (stocks) sheet gains three new trailing columns:
#   H = date, I = legislator_name, J = legislator_id

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$legislatorName = "李鴻鈞"
$legislatorId = 898
$reportDate = "2011-11-17"

# New header cells (row 1)
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Match the bold/bordered header formatting already used by columns B..G
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# Fill the new columns for every existing data row
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Force the date column's data rows to plain text so the date string
# isn't reinterpreted as a date serial number when assigned below.
$ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8)).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = $reportDate
    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}

$excel.CutCopyMode = $false
